# Update F-column ('想去人数' / interest count) values per the commit diff.
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 39
$ws.Range("F5").Value = 7949
$ws.Range("F9").Value = 57
$ws.Range("F10").Value = 6864
$ws.Range("F12").Value = 502
$ws.Range("F13").Value = 473
$ws.Range("F15").Value = 682
$ws.Range("F21").Value = 118
$ws.Range("F22").Value = 11203
$ws.Range("F24").Value = 86
$ws.Range("F25").Value = 2132
$ws.Range("F26").Value = 2885
$ws.Range("F28").Value = 43
$ws.Range("F29").Value = 2539
$ws.Range("F31").Value = 91
$ws.Range("F32").Value = 39
$ws.Range("F34").Value = 2292
$ws.Range("F35").Value = 331
$ws.Range("F36").Value = 1560
$ws.Range("F38").Value = 67
$ws.Range("F39").Value = 5656
$ws.Range("F40").Value = 71
$ws.Range("F41").Value = 1237
$ws.Range("F42").Value = 804
$ws.Range("F43").Value = 150
$ws.Range("F45").Value = 1095
$ws.Range("F46").Value = 1049
$ws.Range("F47").Value = 1478
$ws.Range("F48").Value = 91
$ws.Range("F49").Value = 1120

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F8").Value = 248
$ws.Range("F20").Value = 57
$ws.Range("F22").Value = 2

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 162
$ws.Range("F3").Value = 272

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F4").Value = 162
$ws.Range("F5").Value = 272
$ws.Range("F8").Value = 7949
$ws.Range("F11").Value = 57
$ws.Range("F12").Value = 6864
$ws.Range("F13").Value = 6864
$ws.Range("F15").Value = 502
$ws.Range("F16").Value = 473
$ws.Range("F17").Value = 682
$ws.Range("F21").Value = 248
$ws.Range("F25").Value = 11203
$ws.Range("F27").Value = 86
$ws.Range("F28").Value = 2132
$ws.Range("F29").Value = 2885
$ws.Range("F30").Value = 2539
$ws.Range("F33").Value = 39
$ws.Range("F35").Value = 2292
$ws.Range("F36").Value = 331
$ws.Range("F37").Value = 1560
$ws.Range("F39").Value = 67
$ws.Range("F40").Value = 5656
$ws.Range("F41").Value = 57
$ws.Range("F42").Value = 1237
$ws.Range("F43").Value = 804
$ws.Range("F44").Value = 150
$ws.Range("F46").Value = 1095
$ws.Range("F47").Value = 1049
$ws.Range("F48").Value = 1478
$ws.Range("F49").Value = 91
$ws.Range("F50").Value = 1120
